$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 353
$ws.Range("B353").Value = 6777993
$ws.Range("F353").Value = "Liverpool"
$ws.Range("G353").Value = "Fulham"
$ws.Range("H353").Value = 4
$ws.Range("I353").Value = 3
$ws.Range("J353").Value = "H"
$ws.Range("K353").Value = 1.285
$ws.Range("L353").Value = 6
$ws.Range("M353").Value = 8.5
$ws.Range("N353").Value = 1.25
$ws.Range("O353").Value = 7
$ws.Range("P353").Value = 10
$ws.Range("Q353").Value = -1.75
$ws.Range("R353").Value = 1.97
$ws.Range("S353").Value = 1.96
$ws.Range("T353").Value = 3.25
$ws.Range("U353").Value = 1.87
$ws.Range("V353").Value = 2.03
$ws.Range("W353").Value = 0.25
$ws.Range("X353").Value = -1
$ws.Range("AA353").Value = 0.96
$ws.Range("AB353").Value = 0.8700000000000001
$ws.Range("AC353").Value = -1

# Row 354
$ws.Range("B354").Value = 6777771
$ws.Range("F354").Value = "West Ham"
$ws.Range("G354").Value = "Crystal Palace"
$ws.Range("H354").Value = 1
$ws.Range("I354").Value = 1
$ws.Range("J354").Value = "D"
$ws.Range("K354").Value = 2
$ws.Range("L354").Value = 3.4
$ws.Range("M354").Value = 3.75
$ws.Range("N354").Value = 1.909
$ws.Range("O354").Value = 3.6
$ws.Range("P354").Value = 4
$ws.Range("Q354").Value = -0.5
$ws.Range("R354").Value = 1.95
$ws.Range("S354").Value = 1.95
$ws.Range("T354").Value = 2.25
$ws.Range("U354").Value = 1.88
$ws.Range("V354").Value = 2.02
$ws.Range("W354").Value = -1
$ws.Range("X354").Value = 2.6
$ws.Range("AA354").Value = 0.95
$ws.Range("AB354").Value = -0.5
$ws.Range("AC354").Value = 0.51

# Row 355
$ws.Range("B355").Value = 6777992
$ws.Range("F355").Value = "Chelsea"
$ws.Range("G355").Value = "Brighton"
$ws.Range("H355").Value = 3
$ws.Range("J355").Value = "H"
$ws.Range("K355").Value = 2
$ws.Range("M355").Value = 3.2
$ws.Range("N355").Value = 1.571
$ws.Range("O355").Value = 4.333
$ws.Range("P355").Value = 5.25
$ws.Range("Q355").Value = -1
$ws.Range("R355").Value = 1.97
$ws.Range("S355").Value = 1.93
$ws.Range("T355").Value = 2.75
$ws.Range("U355").Value = 1.84
$ws.Range("V355").Value = 2.06
$ws.Range("W355").Value = 0.571
$ws.Range("X355").Value = -1
$ws.Range("Z355").Value = 0
$ws.Range("AA355").Value = -0
$ws.Range("AB355").Value = 0.8400000000000001

# Row 356
$ws.Range("B356").Value = 6777989
$ws.Range("F356").Value = "Bournemouth"
$ws.Range("G356").Value = "Aston Villa"
$ws.Range("H356").Value = 2
$ws.Range("J356").Value = "D"
$ws.Range("K356").Value = 3.75
$ws.Range("M356").Value = 1.833
$ws.Range("N356").Value = 3
$ws.Range("O356").Value = 4
$ws.Range("P356").Value = 2.1
$ws.Range("Q356").Value = 0.25
$ws.Range("R356").Value = 1.95
$ws.Range("S356").Value = 1.95
$ws.Range("T356").Value = 3.25
$ws.Range("U356").Value = 2.03
$ws.Range("V356").Value = 1.87
$ws.Range("W356").Value = -1
$ws.Range("X356").Value = 3
$ws.Range("Z356").Value = 0.475
$ws.Range("AA356").Value = -0.5
$ws.Range("AB356").Value = 1.03

# Row 361
$ws.Range("B361").Value = 6778000
$ws.Range("F361").Value = "Fulham"
$ws.Range("G361").Value = "Nottm Forest"
$ws.Range("H361").Value = 5
$ws.Range("I361").Value = 0
$ws.Range("J361").Value = "H"
$ws.Range("K361").Value = 2.05
$ws.Range("L361").Value = 3.3
$ws.Range("M361").Value = 3.6
$ws.Range("N361").Value = 2
$ws.Range("O361").Value = 3.4
$ws.Range("P361").Value = 3.75
$ws.Range("Q361").Value = -0.5
$ws.Range("R361").Value = 2.05
$ws.Range("S361").Value = 1.85
$ws.Range("T361").Value = 2.25
$ws.Range("U361").Value = 1.9
$ws.Range("V361").Value = 2
$ws.Range("W361").Value = 1
$ws.Range("Y361").Value = -1
$ws.Range("Z361").Value = 1.05
$ws.Range("AA361").Value = -1
$ws.Range("AB361").Value = 0.8999999999999999
$ws.Range("AC361").Value = -1

# Row 362
$ws.Range("B362").Value = 6778001
$ws.Range("F362").Value = "Sheff Utd"
$ws.Range("G362").Value = "Liverpool"
$ws.Range("K362").Value = 8.5
$ws.Range("L362").Value = 6.5
$ws.Range("M362").Value = 1.25
$ws.Range("N362").Value = 15
$ws.Range("O362").Value = 7.5
$ws.Range("P362").Value = 1.181
$ws.Range("Q362").Value = 2.25
$ws.Range("R362").Value = 1.84
$ws.Range("S362").Value = 2.06
$ws.Range("T362").Value = 3.75
$ws.Range("U362").Value = 2.04
$ws.Range("V362").Value = 1.86
$ws.Range("Y362").Value = 0.181
$ws.Range("Z362").Value = 0.42
$ws.Range("AA362").Value = -0.5
$ws.Range("AB362").Value = -1
$ws.Range("AC362").Value = 0.8600000000000001

# Row 363
$ws.Range("B363").Value = 6778003
$ws.Range("F363").Value = "Crystal Palace"
$ws.Range("G363").Value = "Bournemouth"
$ws.Range("H363").Value = 0
$ws.Range("I363").Value = 2
$ws.Range("J363").Value = "A"
$ws.Range("K363").Value = 2
$ws.Range("L363").Value = 3.4
$ws.Range("N363").Value = 2.25
$ws.Range("O363").Value = 3.5
$ws.Range("P363").Value = 3.1
$ws.Range("Q363").Value = -0.25
$ws.Range("R363").Value = 2.02
$ws.Range("S363").Value = 1.88
$ws.Range("U363").Value = 1.84
$ws.Range("V363").Value = 2.06
$ws.Range("W363").Value = -1
$ws.Range("Y363").Value = 2.1
$ws.Range("Z363").Value = -1
$ws.Range("AA363").Value = 0.8799999999999999
$ws.Range("AB363").Value = -0.5
$ws.Range("AC363").Value = 0.53

# Row 474
$ws.Range("U474").Value = 1.87
$ws.Range("V474").Value = 2.03
